$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values
$ws.Range("C1").Value = "cmsName"
$ws.Range("D1").Value = "cmsBackgroundColor"
$ws.Range("E1").Value = "titleBackgroundColor"
$ws.Range("F1").Value = "titleTextColor"
$ws.Range("G1").Value = "cardBackgroundColor"

# Copy the style used on D1 (existing header style) to the new header cells E1:G1
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Delete row 2 entirely (data row removed)
$ws.Rows("2:2").Delete()
